$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 464
$ws1.Range("F9").Value = 603

# Sheet "全部类型" (All types) - mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 464
$ws4.Range("F9").Value = 603
